$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$values = New-Object 'object[,]' 20,10

$values[0,0] = -19.1658491461764
$values[0,1] = -19.1658491461764
$values[0,2] = -19.1658491461764
$values[0,3] = -19.1658491461764
$values[0,4] = -19.1658491461764
$values[0,5] = -19.1658491461764
$values[0,6] = -19.1658491461764
$values[0,7] = -19.1658491461764
$values[0,8] = -19.1658491461764
$values[0,9] = -19.1658491461764

$values[1,0] = -19.1658491461764
$values[1,1] = -19.1658491461764
$values[1,2] = -19.1658491461764
$values[1,3] = -19.1658491461764
$values[1,4] = -19.1658491461764
$values[1,5] = -19.1658491461764
$values[1,6] = -19.1658491461764
$values[1,7] = 1.246213083088453
$values[1,8] = -19.1658491461764
$values[1,9] = -19.1658491461764

$values[2,0] = -19.1658491461764
$values[2,1] = 4.321925764634718
$values[2,2] = 1.262888246737772
$values[2,3] = -19.1658491461764
$values[2,4] = 3.415440867547367
$values[2,5] = -19.1658491461764
$values[2,6] = 1.474799154089173
$values[2,7] = -19.1658491461764
$values[2,8] = 0.8511182439139189
$values[2,9] = -19.1658491461764

$values[3,0] = -19.1658491461764
$values[3,1] = -19.1658491461764
$values[3,2] = -19.1658491461764
$values[3,3] = -19.1658491461764
$values[3,4] = -19.1658491461764
$values[3,5] = 2.833021173564044
$values[3,6] = -19.1658491461764
$values[3,7] = -19.1658491461764
$values[3,8] = -19.1658491461764
$values[3,9] = -19.1658491461764

$values[4,0] = -19.1658491461764
$values[4,1] = -19.1658491461764
$values[4,2] = -19.1658491461764
$values[4,3] = -19.1658491461764
$values[4,4] = -19.1658491461764
$values[4,5] = -19.1658491461764
$values[4,6] = -19.1658491461764
$values[4,7] = -19.1658491461764
$values[4,8] = -19.1658491461764
$values[4,9] = -19.1658491461764

$values[5,0] = 2.360048685495662
$values[5,1] = -19.1658491461764
$values[5,2] = -19.1658491461764
$values[5,3] = -19.1658491461764
$values[5,4] = -19.1658491461764
$values[5,5] = -19.1658491461764
$values[5,6] = -19.1658491461764
$values[5,7] = -19.1658491461764
$values[5,8] = -19.1658491461764
$values[5,9] = -19.1658491461764

$values[6,0] = -19.1658491461764
$values[6,1] = -19.1658491461764
$values[6,2] = -19.1658491461764
$values[6,3] = 1.883907025514451
$values[6,4] = -19.1658491461764
$values[6,5] = -19.1658491461764
$values[6,6] = -19.1658491461764
$values[6,7] = -19.1658491461764
$values[6,8] = -19.1658491461764
$values[6,9] = -19.1658491461764

$values[7,0] = 3.893953590608627
$values[7,1] = -19.1658491461764
$values[7,2] = -19.1658491461764
$values[7,3] = -19.1658491461764
$values[7,4] = -19.1658491461764
$values[7,5] = -19.1658491461764
$values[7,6] = -19.1658491461764
$values[7,7] = -19.1658491461764
$values[7,8] = -19.1658491461764
$values[7,9] = -19.1658491461764

$values[8,0] = -19.1658491461764
$values[8,1] = -19.1658491461764
$values[8,2] = -19.1658491461764
$values[8,3] = -19.1658491461764
$values[8,4] = -19.1658491461764
$values[8,5] = -19.1658491461764
$values[8,6] = -19.1658491461764
$values[8,7] = 1.736809336258479
$values[8,8] = -19.1658491461764
$values[8,9] = 2.203299268802002

$values[9,0] = -19.1658491461764
$values[9,1] = -19.1658491461764
$values[9,2] = -19.1658491461764
$values[9,3] = 2.919933518942707
$values[9,4] = -19.1658491461764
$values[9,5] = 2.827095066079484
$values[9,6] = -19.1658491461764
$values[9,7] = -19.1658491461764
$values[9,8] = -19.1658491461764
$values[9,9] = 1.977766371364105

$values[10,0] = -19.1658491461764
$values[10,1] = -19.1658491461764
$values[10,2] = -19.1658491461764
$values[10,3] = -19.1658491461764
$values[10,4] = -19.1658491461764
$values[10,5] = -19.1658491461764
$values[10,6] = -19.1658491461764
$values[10,7] = -19.1658491461764
$values[10,8] = -19.1658491461764
$values[10,9] = -19.1658491461764

$values[11,0] = -19.1658491461764
$values[11,1] = -19.1658491461764
$values[11,2] = -19.1658491461764
$values[11,3] = 2.519914693110363
$values[11,4] = -19.1658491461764
$values[11,5] = -19.1658491461764
$values[11,6] = -19.1658491461764
$values[11,7] = -19.1658491461764
$values[11,8] = 1.675147698285796
$values[11,9] = 1.741511733857598

$values[12,0] = -19.1658491461764
$values[12,1] = -19.1658491461764
$values[12,2] = 1.711822171925627
$values[12,3] = -19.1658491461764
$values[12,4] = -19.1658491461764
$values[12,5] = -19.1658491461764
$values[12,6] = -19.1658491461764
$values[12,7] = -19.1658491461764
$values[12,8] = -19.1658491461764
$values[12,9] = 1.967441725085147

$values[13,0] = -19.1658491461764
$values[13,1] = -19.1658491461764
$values[13,2] = 1.843824362860813
$values[13,3] = -19.1658491461764
$values[13,4] = -19.1658491461764
$values[13,5] = -19.1658491461764
$values[13,6] = -19.1658491461764
$values[13,7] = -19.1658491461764
$values[13,8] = -19.1658491461764
$values[13,9] = -19.1658491461764

$values[14,0] = -19.1658491461764
$values[14,1] = -19.1658491461764
$values[14,2] = -19.1658491461764
$values[14,3] = -19.1658491461764
$values[14,4] = -19.1658491461764
$values[14,5] = -19.1658491461764
$values[14,6] = -19.1658491461764
$values[14,7] = -19.1658491461764
$values[14,8] = 1.908329563899977
$values[14,9] = -19.1658491461764

$values[15,0] = -19.1658491461764
$values[15,1] = -19.1658491461764
$values[15,2] = 1.84748238811414
$values[15,3] = -19.1658491461764
$values[15,4] = -19.1658491461764
$values[15,5] = -19.1658491461764
$values[15,6] = 2.070137592125758
$values[15,7] = 2.095417686159475
$values[15,8] = 2.553068717158074
$values[15,9] = -19.1658491461764

$values[16,0] = -19.1658491461764
$values[16,1] = -19.1658491461764
$values[16,2] = -19.1658491461764
$values[16,3] = -19.1658491461764
$values[16,4] = -19.1658491461764
$values[16,5] = -19.1658491461764
$values[16,6] = 2.016984658088967
$values[16,7] = 2.048195406281221
$values[16,8] = 2.427606263967604
$values[16,9] = -19.1658491461764

$values[17,0] = -19.1658491461764
$values[17,1] = -19.1658491461764
$values[17,2] = 1.97774146124569
$values[17,3] = -19.1658491461764
$values[17,4] = -19.1658491461764
$values[17,5] = -19.1658491461764
$values[17,6] = 1.615655450765938
$values[17,7] = 1.811400726454663
$values[17,8] = -19.1658491461764
$values[17,9] = -19.1658491461764

$values[18,0] = -19.1658491461764
$values[18,1] = -19.1658491461764
$values[18,2] = 1.67689787032984
$values[18,3] = -19.1658491461764
$values[18,4] = 3.221926558786009
$values[18,5] = -19.1658491461764
$values[18,6] = 1.64132729598784
$values[18,7] = 1.248300305396875
$values[18,8] = -19.1658491461764
$values[18,9] = 2.070815017665041

$values[19,0] = -19.1658491461764
$values[19,1] = -19.1658491461764
$values[19,2] = -19.1658491461764
$values[19,3] = 1.587675303945201
$values[19,4] = -19.1658491461764
$values[19,5] = 2.530541556786847
$values[19,6] = 1.480190038638628
$values[19,7] = -19.1658491461764
$values[19,8] = -19.1658491461764
$values[19,9] = -19.1658491461764

$ws.Range("B2:K21").Value = $values
